# RPAR_holdings.xlsx update:
#   - bump the "as of" date in the confidential disclaimer (A18) from
#     2021-05-03 to 2021-05-04
#   - refresh the Weight (col D) / Percent Change (col E) figures for the
#     holdings rows (2-15) on Sheet1
#
# The sheet ships with cell protection enabled, so we briefly unprotect it,
# make the edits, then restore protection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# --- Update disclaimer text in A18 (date only) ---
$oldText = $ws.Range("A18").Value
$newText = $oldText -replace "2021-05-03", "2021-05-04"
$ws.Range("A18").Value = $newText

# --- Update Weight (D) / Percent Change (E) values for rows 2-15 ---
$values = @(
    @{ Row = 2;  D = 0.05759327990127471;  E = -0.007494942063637988 }
    @{ Row = 3;  D = 0.02370015902503486;  E = -0.01212358232303479 }
    @{ Row = 4;  D = 0.03155607519328625;  E = -0.01227341389728087 }
    @{ Row = 5;  D = 0.0318877635465893;   E = -0.0003944773175543759 }
    @{ Row = 6;  D = 0.0369425351583034;   E = 0.004731182795698841 }
    @{ Row = 7;  D = 0.01914920944298866;  E = 0.002662149080348453 }
    @{ Row = 8;  D = 0.004550287529444789; E = -0.04103011785246624 }
    @{ Row = 9;  D = 0.006935663052300831; E = 0.004772814051164476 }
    @{ Row = 10; D = 0.0707469409859552;   E = -0.007860752386299685 }
    @{ Row = 11; D = 0.07078666414203941;  E = -0.006734006734006814 }
    @{ Row = 12; D = 0.1467320421542132;   E = 0.006786023678891029 }
    @{ Row = 13; D = 0.3852887939652581;   E = 0.001314521076154707 }
    @{ Row = 14; D = 0.1141305859033112;   E = 0.004019978072846708 }
    @{ Row = 15; D = 0.9999999999999998;   E = -0.0001185074156512256 }
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# --- Restore sheet protection (original had cell protection enabled) ---
$ws.Protect()
